$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new data row (row 11) mirroring the existing rows' structure.
$row = 11

$ws.Cells.Item($row, 1).Value = 7
$ws.Cells.Item($row, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item($row, 3).Value = "Ñuble"

# Column D holds a date value formatted like the cells above it (style index 2 -> YYYY-MM-DD HH:MM:SS)
$ws.Cells.Item($row, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item($row, 4).Value = [datetime]"2021-11-16"

$ws.Cells.Item($row, 5).Value = 16
$ws.Cells.Item($row, 6).Value = 300000000
$ws.Cells.Item($row, 7).Value = "Espárragos"
$ws.Cells.Item($row, 8).Value = "Sin especificar"
$ws.Cells.Item($row, 9).Value = "Primera"
$ws.Cells.Item($row, 10).Value = 400
$ws.Cells.Item($row, 11).Value = 900
$ws.Cells.Item($row, 12).Value = 1000
$ws.Cells.Item($row, 13).Value = 950
$ws.Cells.Item($row, 14).Value = "`$/kilo"
$ws.Cells.Item($row, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item($row, 16).Value = 950
$ws.Cells.Item($row, 17).Value = 1
$ws.Cells.Item($row, 18).Value = "Hortaliza"
